# Update benchmark: 2026-02-11 07:09:57 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "15 TL - 15 TL"
$ws.Range("G3").Value = "39,87 TRY - 79,76 TRY - 797,68 TRY"
$ws.Range("G4").Value = "27,84 TRY - 55,69 TRY - 398,83 TRY"
$ws.Range("G5").Value = "7,97 TRY - 15,96 TRY - 199,41 TRY"
$ws.Range("G6").Value = "8.300,01 TL - 99,71 TL"
$ws.Range("J6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("G8").Value = "19,94 TRY - 39,88 TRY - 398,84 TRY"
$ws.Range("G9").Value = "13,92 TRY - 27,85 TRY - 199,42 TRY"
$ws.Range("G10").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("G11").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"
$ws.Range("G14").Value = "8.300 TL - 7,97 TL"
$ws.Range("J14").Value = "1.554,97 TL - 7.784 TL"
$ws.Range("F17").Value = "%0,5 Asgari Tutar: 361,9 TL Azami Tutar: 361,9 TL / 361,9 TL"
$ws.Range("G17").Value = " Asgari Tutar: 300 TL Azami Tutar: 300 TL"
$ws.Range("F21").Value = "%0,5 Asgari Tutar: 544,76 TL Azami Tutar: 544,76 TL / 3.157,14 TL"
$ws.Range("G21").Value = "%0,16 Asgari Tutar: 300 TL Azami Tutar: "
$ws.Range("F23").Value = "86,67 TL"
$ws.Range("G23").Value = "600 TL"
